$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue 'D2' '29.896.23'
Set-TextValue 'E2' '  +0.07%  '
Set-TextValue 'D3' '1.874.05'
Set-TextValue 'E3' '  -0.72%  '
Set-TextValue 'E4' '  +0.10%  '
Set-TextValue 'D5' '0.7391'
Set-TextValue 'E5' '  -3.80%  '
Set-TextValue 'D6' '242.38'
Set-TextValue 'E6' '  -0.14%  '
Set-TextValue 'E7' '  +0.06%  '
Set-TextValue 'D8' '0.3159'
Set-TextValue 'E8' '  +1.06%  '
Set-TextValue 'D9' '0.07221'
Set-TextValue 'E9' '  +0.71%  '
Set-TextValue 'D10' '24.61'
Set-TextValue 'E10' '  -3.99%  '
Set-TextValue 'D11' '0.08337'
Set-TextValue 'E11' '  -2.82%  '
Set-TextValue 'B12' 'Polygon'
Set-TextValue 'C12' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D12' '0.7494'
Set-TextValue 'E12' '  -1.85%  '
Set-TextValue 'B13' 'WrappedEther'
Set-TextValue 'C13' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D13' '1.888.50'
Set-TextValue 'E13' '  -0.53%  '
Set-TextValue 'D14' '5.394'
Set-TextValue 'E14' '  +0.56%  '
Set-TextValue 'D15' '92.30'
Set-TextValue 'E15' '  -1.41%  '
Set-TextValue 'D16' '29.950.11'
Set-TextValue 'E16' '  +0.36%  '
Set-TextValue 'D17' '6.092'
Set-TextValue 'E17' '  -0.89%  '
Set-TextValue 'D18' '247.66'
Set-TextValue 'E19' '  -1.64%  '
Set-TextValue 'D20' '0.000007838'
Set-TextValue 'E20' '  +0.43%  '
Set-TextValue 'D21' '1.002'
Set-TextValue 'E21' '  +0.33%  '
Set-TextValue 'D22' '2.143.26'
Set-TextValue 'E22' '  +0.66%  '
Set-TextValue 'D23' '8.035'
Set-TextValue 'E23' '  +0.35%  '
Set-TextValue 'D24' '1.000'
Set-TextValue 'E24' '  -0.06%  '
Set-TextValue 'D25' '0.1551'
Set-TextValue 'E25' '  -5.24%  '
Set-TextValue 'D26' '9.265'
Set-TextValue 'E26' '  -1.27%  '
Set-TextValue 'D27' '164.46'
Set-TextValue 'E27' '  +0.92%  '
Set-TextValue 'E28' '  -0.32%  '
Set-TextValue 'D29' '2.027'
Set-TextValue 'E29' '  -0.23%  '
Set-TextValue 'D30' '1.507'
Set-TextValue 'E30' '  +2.65%  '
Set-TextValue 'D31' '4.594'
Set-TextValue 'E31' '  +1.84%  '
Set-TextValue 'E32' '  +0.04%  '
Set-TextValue 'D33' '4.260'
Set-TextValue 'E33' '  +4.07%  '
Set-TextValue 'E34' '  -2.48%  '
Set-TextValue 'E35' '  -0.55%  '
Set-TextValue 'D36' '0.7485'
Set-TextValue 'E36' '  +0.83%  '
Set-TextValue 'D37' '1.001'
Set-TextValue 'E37' '  +0.03%  '
Set-TextValue 'E38' '  -0.05%  '
Set-TextValue 'E39' '  +0.60%  '
Set-TextValue 'D40' '2.752'
Set-TextValue 'E40' '  -1.18%  '
Set-TextValue 'D41' '0.4536'
Set-TextValue 'E41' '  +1.46%  '
Set-TextValue 'B42' 'FraxShare'
Set-TextValue 'C42' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D42' '6.140'
Set-TextValue 'E42' '  +1.15%  '
Set-TextValue 'B43' 'Maker'
Set-TextValue 'C43' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D43' '1.101.93'
Set-TextValue 'E43' '  -0.39%  '
Set-TextValue 'D44' '72.25'
Set-TextValue 'E44' '  -1.19%  '
Set-TextValue 'D45' '0.8625'
Set-TextValue 'E45' '  +1.22%  '
Set-TextValue 'D46' '104.30'
Set-TextValue 'E46' '  +1.68%  '
Set-TextValue 'D47' '1.001'
Set-TextValue 'E47' '  +0.02%  '
Set-TextValue 'D48' '1.859'
Set-TextValue 'D49' '7.605'
Set-TextValue 'E49' '  -0.65%  '
Set-TextValue 'D50' '9.529'
Set-TextValue 'E50' '  -2.34%  '
Set-TextValue 'D51' '2.041.03'
Set-TextValue 'E51' '  +0.25%  '
